$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("snapshot")

# Updated scraped_at timestamps (column K) for rows 2-43 on the "snapshot" sheet.
# These mirror a fresh scrape re-run; all other cell data is unchanged.
$ws.Range("K2").Value = "2025-11-17T23:21:26.318468+00:00"
$ws.Range("K3").Value = "2025-11-17T23:21:26.318504+00:00"
$ws.Range("K4").Value = "2025-11-17T23:21:26.318524+00:00"
$ws.Range("K5").Value = "2025-11-17T23:21:28.842913+00:00"
$ws.Range("K6").Value = "2025-11-17T23:21:28.842956+00:00"
$ws.Range("K7").Value = "2025-11-17T23:21:28.842968+00:00"
$ws.Range("K8").Value = "2025-11-17T23:21:31.231401+00:00"
$ws.Range("K9").Value = "2025-11-17T23:21:33.674152+00:00"
$ws.Range("K10").Value = "2025-11-17T23:21:35.774131+00:00"
$ws.Range("K11").Value = "2025-11-17T23:21:35.774160+00:00"
$ws.Range("K12").Value = "2025-11-17T23:21:40.753322+00:00"
$ws.Range("K13").Value = "2025-11-17T23:21:43.179370+00:00"
$ws.Range("K14").Value = "2025-11-17T23:21:45.691118+00:00"
$ws.Range("K15").Value = "2025-11-17T23:21:45.691139+00:00"
$ws.Range("K16").Value = "2025-11-17T23:21:45.691147+00:00"
$ws.Range("K17").Value = "2025-11-17T23:21:48.199767+00:00"
$ws.Range("K18").Value = "2025-11-17T23:21:50.648765+00:00"
$ws.Range("K19").Value = "2025-11-17T23:21:50.648782+00:00"
$ws.Range("K20").Value = "2025-11-17T23:21:52.821423+00:00"
$ws.Range("K21").Value = "2025-11-17T23:21:54.918571+00:00"
$ws.Range("K22").Value = "2025-11-17T23:21:54.918600+00:00"
$ws.Range("K23").Value = "2025-11-17T23:21:54.918618+00:00"
$ws.Range("K24").Value = "2025-11-17T23:21:54.918634+00:00"
$ws.Range("K25").Value = "2025-11-17T23:21:57.371555+00:00"
$ws.Range("K26").Value = "2025-11-17T23:21:57.371579+00:00"
$ws.Range("K27").Value = "2025-11-17T23:21:59.929884+00:00"
$ws.Range("K28").Value = "2025-11-17T23:21:59.929913+00:00"
$ws.Range("K29").Value = "2025-11-17T23:21:59.929933+00:00"
$ws.Range("K30").Value = "2025-11-17T23:22:02.388274+00:00"
$ws.Range("K31").Value = "2025-11-17T23:22:02.388316+00:00"
$ws.Range("K32").Value = "2025-11-17T23:22:04.818936+00:00"
$ws.Range("K33").Value = "2025-11-17T23:22:04.818971+00:00"
$ws.Range("K34").Value = "2025-11-17T23:22:04.818991+00:00"
$ws.Range("K35").Value = "2025-11-17T23:22:04.819009+00:00"
$ws.Range("K36").Value = "2025-11-17T23:22:04.819024+00:00"
$ws.Range("K37").Value = "2025-11-17T23:22:07.484526+00:00"
$ws.Range("K38").Value = "2025-11-17T23:22:07.484556+00:00"
$ws.Range("K39").Value = "2025-11-17T23:22:12.529691+00:00"
$ws.Range("K40").Value = "2025-11-17T23:22:12.529722+00:00"
$ws.Range("K41").Value = "2025-11-17T23:22:12.529741+00:00"
$ws.Range("K42").Value = "2025-11-17T23:22:12.529758+00:00"
$ws.Range("K43").Value = "2025-11-17T23:22:14.922456+00:00"
